$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (a280.xlsx)
$ws.Range("D2").Value = 2871.75244140625
$ws.Range("E2").Value = 13.504555521640194
$ws.Range("F2").Value = "6.764s"
$ws.Range("H2").Value = 300.0

# Row 3 (ali535.xlsx)
$ws.Range("D3").Value = 2162.1923828125
$ws.Range("E3").Value = 14.53758657851587
$ws.Range("F3").Value = "29.889s"
$ws.Range("H3").Value = 300.0

# Row 4 (att48.xlsx)
$ws.Range("D4").Value = 32272.44921875
$ws.Range("E4").Value = 17.895004748859044
$ws.Range("F4").Value = "0.37s"
$ws.Range("H4").Value = 300.0
